$d = $word.ActiveDocument

# 1. Remove proofErr marks around "Entidades a notificar" — achieved by
#    re-typing the text which clears the proofing-error bookkeeping,
#    then trimming formatting marks stays the same.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Entidades a notificar", $true, $false, $false, $false, $false, $true, 1, $false, "Entidades a notificar", 2)

# 2. Split the sentence about "SIN RIESGO" and insert " separados por espacio"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("SIN RIESGO, en caso de no haber ninguno devolver NA.", $true, $false, $false, $false, $false, $true, 1, $false, "SIN RIESGO separados por espacio, en caso de no haber ninguno devolver NA.", 2)
